# Updated template file with versioning in README
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReadMe")

# Insert a new row above the existing "Version updates" entries (row 10)
# so the "V0.0 / Original version" and "V0.1 / ..." rows shift down.
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "This is version V0.1"

$ws.Range("A14").Select()
